$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 (A4:DK4) into row 5 (A5:DK5) so that the
# blank "inlineStr" placeholder cells (D, G, AG, AX) and the header-style
# cell A5 get created/styled the same way as the existing rows.
$ws.Range("A4:DK4").Copy()
$ws.Range("A5:DK5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row label
$ws.Range("A5").Value = "2021年"

# Numeric data for 2021
$ws.Range("B5").Value = 15.8
$ws.Range("C5").Value = 30.3
$ws.Range("E5").Value = 65.40000000000001
$ws.Range("F5").Value = 28.9
$ws.Range("H5").Value = 21.3
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 5.2
$ws.Range("K5").Value = -12.9
$ws.Range("L5").Value = 7.8
$ws.Range("M5").Value = 42.5
$ws.Range("N5").Value = 0.1
$ws.Range("O5").Value = 8.1
$ws.Range("P5").Value = 8.4
$ws.Range("Q5").Value = -5
$ws.Range("R5").Value = 84.09999999999999
$ws.Range("S5").Value = -37.3
$ws.Range("T5").Value = 6.6
$ws.Range("U5").Value = 10.6
$ws.Range("V5").Value = -2.1
$ws.Range("W5").Value = 24.6
$ws.Range("X5").Value = 23.1
$ws.Range("Y5").Value = 18.5
$ws.Range("Z5").Value = 38.5
$ws.Range("AA5").Value = 13.8
$ws.Range("AB5").Value = 10.5
$ws.Range("AC5").Value = 16
$ws.Range("AD5").Value = 12.5
$ws.Range("AE5").Value = 28.2
$ws.Range("AF5").Value = 19.5
$ws.Range("AH5").Value = -49.7
$ws.Range("AI5").Value = 38.5
$ws.Range("AJ5").Value = -22.6
$ws.Range("AK5").Value = 30.3
$ws.Range("AL5").Value = 11.1
$ws.Range("AM5").Value = 4.2
$ws.Range("AN5").Value = -8.4
$ws.Range("AO5").Value = -8.9
$ws.Range("AP5").Value = -33.2
$ws.Range("AQ5").Value = 19.6
$ws.Range("AR5").Value = 57.9
$ws.Range("AS5").Value = -17.8
$ws.Range("AT5").Value = 26.9
$ws.Range("AU5").Value = 82.5
$ws.Range("AV5").Value = 6.2
$ws.Range("AW5").Value = 122.6
$ws.Range("AY5").Value = 3
$ws.Range("AZ5").Value = 41.9
$ws.Range("BA5").Value = 6.2
$ws.Range("BB5").Value = 3.4
$ws.Range("BC5").Value = 12.7
$ws.Range("BD5").Value = 93.3
$ws.Range("BE5").Value = 22.7
$ws.Range("BF5").Value = 22.5
$ws.Range("BG5").Value = 15.8
$ws.Range("BH5").Value = -8.800000000000001
$ws.Range("BI5").Value = -9.5
$ws.Range("BJ5").Value = 19.5
$ws.Range("BK5").Value = 81.2
$ws.Range("BL5").Value = 10.3
$ws.Range("BM5").Value = 35.4
$ws.Range("BN5").Value = 10.9
$ws.Range("BO5").Value = 17.1
$ws.Range("BP5").Value = 4.4
$ws.Range("BQ5").Value = 26.3
$ws.Range("BR5").Value = 22.8
$ws.Range("BS5").Value = -8.1
$ws.Range("BT5").Value = 19.9
$ws.Range("BU5").Value = 27
$ws.Range("BV5").Value = 15.4
$ws.Range("BW5").Value = 20
$ws.Range("BX5").Value = 31.3
$ws.Range("BY5").Value = 23.3
$ws.Range("BZ5").Value = 0.6
$ws.Range("CA5").Value = 19.5
$ws.Range("CB5").Value = 48.8
$ws.Range("CC5").Value = 15.6
$ws.Range("CD5").Value = -16.2
$ws.Range("CE5").Value = 2.1
$ws.Range("CF5").Value = 17.6
$ws.Range("CG5").Value = 19.7
$ws.Range("CH5").Value = 24.2
$ws.Range("CI5").Value = 28.1
$ws.Range("CJ5").Value = -6.2
$ws.Range("CK5").Value = 17
$ws.Range("CL5").Value = 9.199999999999999
$ws.Range("CM5").Value = -44.1
$ws.Range("CN5").Value = 36.5
$ws.Range("CO5").Value = 12
$ws.Range("CP5").Value = 43.3
$ws.Range("CQ5").Value = 62.4
$ws.Range("CR5").Value = -28
$ws.Range("CS5").Value = 26.8
$ws.Range("CT5").Value = 18.1
$ws.Range("CU5").Value = 26.9
$ws.Range("CV5").Value = 40.9
$ws.Range("CW5").Value = 38.9
$ws.Range("CX5").Value = 25.5
$ws.Range("CY5").Value = 27.4
$ws.Range("CZ5").Value = 57.7
$ws.Range("DA5").Value = 19
$ws.Range("DB5").Value = -3.6
$ws.Range("DC5").Value = 27.7
$ws.Range("DD5").Value = 47
$ws.Range("DE5").Value = 5.5
$ws.Range("DF5").Value = 21.8
$ws.Range("DG5").Value = 27
$ws.Range("DH5").Value = 15.6
$ws.Range("DI5").Value = 9.5
$ws.Range("DJ5").Value = 40.5
$ws.Range("DK5").Value = 35.7
